$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (Grau de entendimento / Grau de representação table)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2

$ws.Range("F4").Value = 2

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 3

# Update the selected cell shown in the sheet view
$ws.Range("F9").Select()
